# The commit swaps the colour scheme carried by the deck's two theme
# parts: ppt/theme/theme1.xml (the slide master's theme, currently the
# "Integral" / "Red Violet" palette) ends up holding the "Office Theme"
# palette that used to live in ppt/theme/theme2.xml (the notes master's
# theme), and vice versa. The <a:fontScheme>/<a:fmtScheme> blocks are
# byte-for-byte identical between the two themes already, so the only
# observable difference is the <a:clrScheme> RGB values (the cosmetic
# name="..." attributes are not settable through the PowerPoint object
# model, in real PowerPoint or here).
#
# PowerPoint's automation surface exposes the *active* (slide-master)
# theme's colours through Slide.ThemeColorScheme - each of the 12 slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) maps straight onto the
# corresponding element inside ppt/theme/theme1.xml's <a:clrScheme>, and
# writing its .RGB property edits that XML in place.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeRGB {
    param($item, [string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $item.RGB = $r + ($g * 256) + ($b * 65536)
}

# Target values: the "Office Theme" colour scheme that currently lives
# in ppt/theme/theme2.xml, applied here to the slide master's theme
# (ppt/theme/theme1.xml), slot-for-slot in clrScheme order.
Set-ThemeRGB $tcs.Item(1)  "000000"   # dk1
Set-ThemeRGB $tcs.Item(2)  "FFFFFF"   # lt1
Set-ThemeRGB $tcs.Item(3)  "44546A"   # dk2
Set-ThemeRGB $tcs.Item(4)  "E7E6E6"   # lt2
Set-ThemeRGB $tcs.Item(5)  "5B9BD5"   # accent1
Set-ThemeRGB $tcs.Item(6)  "ED7D31"   # accent2
Set-ThemeRGB $tcs.Item(7)  "A5A5A5"   # accent3
Set-ThemeRGB $tcs.Item(8)  "FFC000"   # accent4
Set-ThemeRGB $tcs.Item(9)  "4472C4"   # accent5
Set-ThemeRGB $tcs.Item(10) "70AD47"   # accent6
Set-ThemeRGB $tcs.Item(11) "0563C1"   # hlink
Set-ThemeRGB $tcs.Item(12) "954F72"   # folHlink
